$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Requirement T3.1.2: MoSCoW marker "M" -> "S"
# ---------------------------------------------------------------------------
$tbl = $d.Tables.Item(1)
foreach ($row in $tbl.Rows) {
    $idCell = $row.Cells.Item(1)
    if ($idCell.Range.Text -like "T3.1.2*") {
        $mCell = $row.Cells.Item(2)
        $mCell.Range.Text = "S"
    }
}

# ---------------------------------------------------------------------------
# 2. "drv8838 will be utilized" -> "drv8838 should be utilized"
#    ("will" becomes "should", emitted as its own run so that "should" is
#    bold/bCs like its neighbours but is a distinct <w:r>.)
# ---------------------------------------------------------------------------
$hit = $d.Content
$hit.Find.Execute("drv8838 will be", $false)
$segStart = $hit.Start
$segText = $hit.Text
$willOffset = $segText.IndexOf("will")
$willStart = $segStart + $willOffset
$willEnd = $willStart + 4

$willRange = $d.Range($willStart, $willEnd)
$willRange.Text = "should"

# Reacquire the "should" range (text length grew by 2 chars) and force it
# onto its own run boundary by toggling a character property off/on.
$shouldEnd = $willStart + 6
$shouldRange = $d.Range($willStart, $shouldEnd)
$shouldRange.Bold = $false
$shouldRange.Bold = $true

# Likewise re-isolate the pre-existing "utilized" run (immediately after
# " be ") so it stays separate from the newly split " be " run, searching
# only from just after "should" onward so the earlier "utilized" occurring
# elsewhere in the table (T4.1 row) can never match.
$utilizedRange = $d.Range($shouldEnd, $d.Content.End)
$utilizedRange.Find.Execute("utilized", $false)
$utilizedRange.Bold = $false
$utilizedRange.Bold = $true
